$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Category" column (H), populated in an order that first introduces
# Category, then CategoryA, then CategoryC, then CategoryB (matches how the
# unique strings ended up ordered in the shared string table).
$ws.Range("H1").Value = "Category"
$ws.Range("H2").Value = "CategoryA"
$ws.Range("H8").Value = "CategoryC"
$ws.Range("H3").Value = "CategoryB"
$ws.Range("H4").Value = "CategoryA"
$ws.Range("H5").Value = "CategoryB"
$ws.Range("H6").Value = "CategoryA"
$ws.Range("H7").Value = "CategoryA"
$ws.Range("H9").Value = "CategoryA"
$ws.Range("H10").Value = "CategoryA"
$ws.Range("H11").Value = "CategoryC"

$ws.Range("H12").Select()
